# Update formant frequency defaults (E column) based on revised source data
# Reference: 10.1121/1.426686
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("E89").Value = 816.5
$ws.Range("E90").Value = 2148.5
$ws.Range("E91").Value = 3077.5
$ws.Range("E92").Value = 729.5
$ws.Range("E93").Value = 2135.5
$ws.Range("E94").Value = 2930
$ws.Range("E95").Value = 704.5
$ws.Range("E96").Value = 2112.5
$ws.Range("E97").Value = 2935.5
$ws.Range("E98").Value = 691.5
$ws.Range("E99").Value = 2106
$ws.Range("E100").Value = 2919
$ws.Range("E101").Value = 696.5
$ws.Range("E102").Value = 2074.5
$ws.Range("E103").Value = 2820
$ws.Range("E104").Value = 697
$ws.Range("E105").Value = 2060.5
$ws.Range("E106").Value = 2822.5
$ws.Range("E107").Value = 650
$ws.Range("E108").Value = 1990
$ws.Range("E109").Value = 2725
$ws.Range("E110").Value = 624.5
$ws.Range("E111").Value = 1919
$ws.Range("E112").Value = 2678
$ws.Range("E113").Value = 576.5
$ws.Range("E114").Value = 1777.5
$ws.Range("E115").Value = 2568.5
$ws.Range("E116").Value = 597
$ws.Range("E117").Value = 1861.5
$ws.Range("E118").Value = 2607
$ws.Range("E119").Value = 520.5
$ws.Range("E120").Value = 1663
$ws.Range("E121").Value = 2359.5
$ws.Range("E122").Value = 387.5
$ws.Range("E123").Value = 1634.5
$ws.Range("E124").Value = 2392
$ws.Range("E125").Value = 501
$ws.Range("E126").Value = 1604
$ws.Range("E127").Value = 2398
$ws.Range("E128").Value = 510
$ws.Range("E129").Value = 1585
$ws.Range("E130").Value = 2337.5
$ws.Range("E131").Value = 507.5
$ws.Range("E132").Value = 1597.5
$ws.Range("E133").Value = 2271
$ws.Range("E134").Value = 845
$ws.Range("E135").Value = 2178
$ws.Range("E136").Value = 3012.5
$ws.Range("E137").Value = 798
$ws.Range("E138").Value = 2122
$ws.Range("E139").Value = 3108.5
$ws.Range("E140").Value = 767
$ws.Range("E141").Value = 2152.5
$ws.Range("E142").Value = 3031.5
$ws.Range("E143").Value = 767
$ws.Range("E144").Value = 2135.5
$ws.Range("E145").Value = 2992.5
$ws.Range("E146").Value = 759
$ws.Range("E147").Value = 2155.5
$ws.Range("E148").Value = 2962
$ws.Range("E149").Value = 754.5
$ws.Range("E150").Value = 2116.5
$ws.Range("E151").Value = 2885
$ws.Range("E152").Value = 723.5
$ws.Range("E153").Value = 2095
$ws.Range("E154").Value = 2814.5
$ws.Range("E155").Value = 689
$ws.Range("E156").Value = 2048
$ws.Range("E157").Value = 2811.5
$ws.Range("E158").Value = 692.5
$ws.Range("E159").Value = 2036.5
$ws.Range("E160").Value = 2776
$ws.Range("E161").Value = 654
$ws.Range("E162").Value = 1942.5
$ws.Range("E163").Value = 2668.5
$ws.Range("E164").Value = 639
$ws.Range("E165").Value = 1915
$ws.Range("E166").Value = 2593
$ws.Range("E167").Value = 637
$ws.Range("E168").Value = 1967.5
$ws.Range("E169").Value = 2648.5
$ws.Range("E170").Value = 656
$ws.Range("E171").Value = 1941.5
$ws.Range("E172").Value = 2642.5
$ws.Range("E173").Value = 675
$ws.Range("E174").Value = 1983
$ws.Range("E175").Value = 2681.5
$ws.Range("E176").Value = 627
$ws.Range("E177").Value = 1918
$ws.Range("E178").Value = 2587.5

# Restore sheet selection / scroll position to match the editor
[void]$ws.Activate()
[void]$ws.Range("B124").Select()
